# Update predicted Output Values (column C) on Sheet1 per the
# "modified predictor, finished base weight function" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 575.9299999999999
$ws.Range("C3").Value = 587.0599999999999
$ws.Range("C4").Value = 557.61
$ws.Range("C5").Value = 581.3
$ws.Range("C6").Value = 581.3
